$d = $word.ActiveDocument

# Find the "Requisitos" paragraph that reads
# "LOM3202: Circuitos Elétricos (Indicação de Conjunto)". Right after it the
# page used to have a blank paragraph followed by a "Ver no Jupiter ..."
# paragraph and a "© 2020 ..." footer paragraph; the site rebuild removed
# those three paragraphs, leaving the "LOM3202..." paragraph immediately
# followed by the blank paragraph that used to sit right before the page
# break at the end of the document.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOM3202: Circuitos Elétricos (Indicação de Conjunto)*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $p1 = $target.Next()
    $p2 = $p1.Next()
    $p3 = $p2.Next()

    if (($p2.Range.Text -like "*Ver no Jupiter*") -and ($p3.Range.Text -like "*Contact: luizeleno@usp.br*")) {
        $r = $d.Range($target.Range.End, $p3.Range.End)
        $r.Delete()
    }
}
